# Update workbook to share with Christopher:
#  - remove the old ad-hoc sample/treatment/description table from SraRunTable (H33:J45)
#  - add two new sheets: "bisulfite_treats" and "rna_treats" with cleaned-up metadata
#  - make "rna_treats" the active tab

$wb = $excel.ActiveWorkbook
$sra = $wb.Worksheets.Item("SraRunTable")

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets (bisulfite_treats first, then rna_treats),
#    right after SraRunTable, then reorder so rna_treats comes first.
# ---------------------------------------------------------------------------
$bisulfite = $wb.Worksheets.Add($null, $sra)
$bisulfite.Name = "bisulfite_treats"

$rna = $wb.Worksheets.Add($null, $bisulfite)
$rna.Name = "rna_treats"

$rna.Move($bisulfite)

# ---------------------------------------------------------------------------
# 2. Populate rna_treats: run / treat / treatInfo for the 12 RNA-Seq runs
# ---------------------------------------------------------------------------
$rna = $wb.Worksheets.Item("rna_treats")

$rna.Cells.Item(1,1).Value = "run"
$rna.Cells.Item(1,2).Value = "treat"
$rna.Cells.Item(1,3).Value = "treatInfo"

$rnaRows = @(
    @("SRR3139744","control","final instar larva (worker)"),
    @("SRR3139743","control","final instar larva (worker)"),
    @("SRR3139742","control","final instar larva (worker)"),
    @("SRR3139741","control","final instar larva (worker)"),
    @("SRR3139740","control","final instar larva (worker)"),
    @("SRR3139739","control","final instar larva (worker)"),
    @("SRR3139738","treated","adult alate"),
    @("SRR3139737","treated","adult alate"),
    @("SRR3139736","treated","adult alate"),
    @("SRR3139735","treated","adult alate"),
    @("SRR3139734","treated","adult alate"),
    @("SRR3139733","treated","adult alate")
)

$r = 2
foreach ($row in $rnaRows) {
    $rna.Cells.Item($r,1).Value = $row[0]
    $rna.Cells.Item($r,2).Value = $row[1]
    $rna.Cells.Item($r,3).Value = $row[2]
    $rna.Cells.Item($r,1).Font.Color = 0
    $rna.Cells.Item($r,3).Font.Color = 0
    $r = $r + 1
}

$rna.Range("A1:A1").ColumnWidth = 10.666666666666666
$rna.Range("C1:C1").ColumnWidth = 21.498697916666668

$rna.Range("A1:C13").Select()

# ---------------------------------------------------------------------------
# 3. Populate bisulfite_treats: run / treatInfo / id / treat for the 8
#    Bisulfite-Seq runs (bismark coverage file names)
# ---------------------------------------------------------------------------
$bisulfite = $wb.Worksheets.Item("bisulfite_treats")

$bisulfite.Cells.Item(1,1).Value = "run"
$bisulfite.Cells.Item(1,2).Value = "treatInfo"
$bisulfite.Cells.Item(1,3).Value = "id"
$bisulfite.Cells.Item(1,4).Value = "treat"

$bisulfiteRows = @(
    @("SRR3139752_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","final instar larva (worker)","larva1",0),
    @("SRR3139751_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","final instar larva (worker)","larva2",0),
    @("SRR3139750_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","final instar larva (worker)","larva3",0),
    @("SRR3139749_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","final instar larva (worker)","larva4",0),
    @("SRR3139748_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","adult alate","adult1",1),
    @("SRR3139747_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","adult alate","adult2",1),
    @("SRR3139746_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","adult alate","adult3",1),
    @("SRR3139745_1.trim_bismark_bt2_pe.deduplicated.bismark.cov","adult alate","adult4",1)
)

$r = 2
foreach ($row in $bisulfiteRows) {
    $bisulfite.Cells.Item($r,1).Value = $row[0]
    $bisulfite.Cells.Item($r,2).Value = $row[1]
    $bisulfite.Cells.Item($r,3).Value = $row[2]
    $bisulfite.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

$bisulfite.Range("A1:A1").ColumnWidth = 10.666666666666666
$bisulfite.Range("B1:B1").ColumnWidth = 21.498697916666668

$bisulfite.Range("C18").Select()

# ---------------------------------------------------------------------------
# 4. Remove the old sample/treatment/description table from SraRunTable
#    (rows 33-45, columns H:J) - this data now lives in rna_treats
# ---------------------------------------------------------------------------
$sra = $wb.Worksheets.Item("SraRunTable")
$sra.Range("A33:AL45").EntireRow.Delete()
$sra.Range("J2:J9").Select()

# ---------------------------------------------------------------------------
# 5. Make rna_treats the active tab
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("rna_treats").Activate()
